$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update violating data points (zone mixture rule violations)
$ws.Range("B18").Value = 3
$ws.Range("B20").Value = -2.6
$ws.Range("B21").Value = 3
$ws.Range("B23").Value = 2.6
$ws.Range("B24").Value = -2.7
$ws.Range("B54").Value = -0.73

# Update the view's top-left cell and selection to match the new scroll
# position (author scrolled down to review the newly flagged points)
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B55").Select()
